# Scheduled market-data refresh: update crafting profit figures
# (currentAveragePrice / LevePrice / LeveProfit columns H-N) on every
# Leve worksheet with freshly pulled numbers.
$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 926.5454999999999
$ws.Range("I12").Value = 926.5454999999999
$ws.Range("K12").Value = 926.5454999999999
$ws.Range("M12").Value = -756.5454999999999
$ws.Range("H17").Value = 3281.4285
$ws.Range("J17").Value = 3985
$ws.Range("L17").Value = 11955
$ws.Range("N17").Value = -12291
$ws.Range("H19").Value = 382.27274
$ws.Range("J19").Value = 133.5
$ws.Range("L19").Value = 133.5
$ws.Range("N19").Value = -483.5
$ws.Range("H74").Value = 4417.8335
$ws.Range("I74").Value = 3003
$ws.Range("J74").Value = 4700.8
$ws.Range("K74").Value = 3003
$ws.Range("L74").Value = 4700.8
$ws.Range("M74").Value = -2067
$ws.Range("N74").Value = -6572.8
$ws.Range("H77").Value = 4417.8335
$ws.Range("I77").Value = 3003
$ws.Range("J77").Value = 4700.8
$ws.Range("K77").Value = 15015
$ws.Range("L77").Value = 23504
$ws.Range("M77").Value = -10335
$ws.Range("N77").Value = -32864
$ws.Range("H100").Value = 2702.75
$ws.Range("I100").Value = 2562.5
$ws.Range("K100").Value = 2562.5
$ws.Range("M100").Value = -2021.5
$ws.Range("H111").Value = 3622.8333
$ws.Range("I111").Value = 3747.4
$ws.Range("J111").Value = 3000
$ws.Range("K111").Value = 11242.2
$ws.Range("L111").Value = 9000
$ws.Range("M111").Value = -8175.200000000001
$ws.Range("N111").Value = -15134
$ws.Range("H113").Value = 5419.4
$ws.Range("I113").Value = 4899.25
$ws.Range("K113").Value = 4899.25
$ws.Range("M113").Value = -1645.25

# --- Worksheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 11473.5
$ws.Range("J29").Value = 11473.5
$ws.Range("L29").Value = 11473.5
$ws.Range("N29").Value = -12089.5
$ws.Range("H45").Value = 2197.4
$ws.Range("I45").Value = 2197.4
$ws.Range("K45").Value = 2197.4
$ws.Range("M45").Value = -1820.4
$ws.Range("H110").Value = 62501076
$ws.Range("I110").Value = 649.5
$ws.Range("K110").Value = 649.5
$ws.Range("M110").Value = 1395.5
$ws.Range("H124").Value = 83142.664
$ws.Range("J124").Value = 83142.664
$ws.Range("L124").Value = 83142.664
$ws.Range("N124").Value = -92962.664

# --- Worksheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 100698.5
$ws.Range("I107").Value = 100698.5
$ws.Range("K107").Value = 100698.5
$ws.Range("M107").Value = -98778.5

# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 585.59375
$ws.Range("I22").Value = 756
$ws.Range("J22").Value = 435.2353
$ws.Range("K22").Value = 756
$ws.Range("L22").Value = 435.2353
$ws.Range("M22").Value = -406
$ws.Range("N22").Value = -1135.2353
$ws.Range("H42").Value = 34000
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 34000
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 34000
$ws.Range("N42").Value = -35186
$ws.Range("M42").ClearContents()

# --- Worksheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 81581.75999999999
$ws.Range("I4").Value = 1495
$ws.Range("J4").Value = 335189.84
$ws.Range("K4").Value = 4485
$ws.Range("L4").Value = 1005569.52
$ws.Range("M4").Value = -4373
$ws.Range("N4").Value = -1005793.52
$ws.Range("H34").Value = 19016.445
$ws.Range("J34").Value = 21374.875
$ws.Range("L34").Value = 64124.625
$ws.Range("N34").Value = -64292.625
$ws.Range("H64").Value = 2000
$ws.Range("J64").Value = 2000
$ws.Range("L64").Value = 6000
$ws.Range("N64").Value = -6540
$ws.Range("H67").Value = 2000
$ws.Range("J67").Value = 2000
$ws.Range("L67").Value = 6000
$ws.Range("N67").Value = -7872
$ws.Range("H68").Value = 2026.8
$ws.Range("I68").Value = 2124.3333
$ws.Range("K68").Value = 6372.999899999999
$ws.Range("M68").Value = -5561.999899999999
$ws.Range("H71").Value = 2026.8
$ws.Range("I71").Value = 2124.3333
$ws.Range("K71").Value = 19118.9997
$ws.Range("M71").Value = -15062.9997
$ws.Range("H121").Value = 458.93332
$ws.Range("I121").Value = 185
$ws.Range("J121").Value = 772
$ws.Range("K121").Value = 555
$ws.Range("L121").Value = 2316
$ws.Range("M121").Value = 755
$ws.Range("N121").Value = -4936
$ws.Range("H134").Value = 4987.1665
$ws.Range("I134").Value = 4987.1665
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 14961.4995
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -9891.499500000002
$ws.Range("H140").Value = 2999.3333
$ws.Range("I140").Value = 2999.3333
$ws.Range("K140").Value = 8997.999899999999
$ws.Range("M140").Value = -3817.999899999999
$ws.Range("H141").Value = 2559.8
$ws.Range("I141").Value = 1966.3334
$ws.Range("J141").Value = 3450
$ws.Range("K141").Value = 5899.0002
$ws.Range("L141").Value = 10350
$ws.Range("M141").Value = -719.0002000000004
$ws.Range("N141").Value = -20710
$ws.Range("N134").ClearContents()

# --- Worksheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 20999
$ws.Range("J40").Value = 20999
$ws.Range("L40").Value = 20999
$ws.Range("N40").Value = -21301
$ws.Range("H47").Value = 29997.5
$ws.Range("J47").Value = 29997.5
$ws.Range("L47").Value = 29997.5
$ws.Range("N47").Value = -31133.5
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("H113").Value = 1124.8334
$ws.Range("I113").Value = 1187.25
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 1187.25
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 982.75
$ws.Range("N113").Value = -5340
$ws.Range("H122").Value = 4245.107
$ws.Range("I122").Value = 2648.8333
$ws.Range("K122").Value = 7946.499899999999
$ws.Range("M122").Value = -5496.499899999999
$ws.Range("H141").Value = 100000
$ws.Range("J141").Value = 100000
$ws.Range("L141").Value = 100000
$ws.Range("N141").Value = -110360
$ws.Range("N75").ClearContents()
$ws.Range("N78").ClearContents()

# --- Worksheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 782.53845
$ws.Range("I22").Value = 672
$ws.Range("J22").Value = 959.4
$ws.Range("K22").Value = 672
$ws.Range("L22").Value = 959.4
$ws.Range("M22").Value = -377
$ws.Range("N22").Value = -1549.4
$ws.Range("H27").Value = 782.53845
$ws.Range("I27").Value = 672
$ws.Range("J27").Value = 959.4
$ws.Range("K27").Value = 672
$ws.Range("L27").Value = 959.4
$ws.Range("M27").Value = -565
$ws.Range("N27").Value = -1173.4
$ws.Range("H35").Value = 12648.625
$ws.Range("I35").Value = 859.5
$ws.Range("J35").Value = 16578.334
$ws.Range("K35").Value = 859.5
$ws.Range("L35").Value = 16578.334
$ws.Range("M35").Value = -523.5
$ws.Range("N35").Value = -17250.334
$ws.Range("H103").Value = 21200.666
$ws.Range("J103").Value = 21200.666
$ws.Range("L103").Value = 21200.666
$ws.Range("N103").Value = -23544.666
$ws.Range("H122").Value = 3452.1667
$ws.Range("I122").Value = 3178.25
$ws.Range("K122").Value = 9534.75
$ws.Range("M122").Value = -7084.75

# --- Worksheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 24000
$ws.Range("J38").Value = 22000
$ws.Range("L38").Value = 22000
$ws.Range("N38").Value = -22946
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("H107").Value = 468.8125
$ws.Range("I107").Value = 333.41666
$ws.Range("J107").Value = 875
$ws.Range("K107").Value = 1000.24998
$ws.Range("L107").Value = 2625
$ws.Range("M107").Value = 919.7500200000001
$ws.Range("N107").Value = -6465
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()
